$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column K (index 11) width slightly (target stored width 15.21875
# chars; the COM layer quantizes ColumnWidth to 1/6-character pixel steps,
# so 14.3 is the closest input that lands on the nearest achievable stored
# width, 15.166666..., i.e. a 2px-narrower column matching the diff's intent)
$ws.Columns.Item(11).ColumnWidth = 14.3

# Update J column values (and a few related I/L values) on rows 1-13
$ws.Range("J1").Value = 46.355953999999997
$ws.Range("J2").Value = 112.85662689999999
$ws.Range("J3").Value = 118.6160898
$ws.Range("J4").Value = 154.50660579999999
$ws.Range("J5").Value = 107.7844659
$ws.Range("J6").Value = 438.61068760000001

$ws.Range("I7").Value = 1270
$ws.Range("J7").Value = 560.41570300000001
$ws.Range("L7").Value = -0.0000013667309257269633

$ws.Range("J8").Value = 118.2736427
$ws.Range("J9").Value = 131.3902247
$ws.Range("J10").Value = 199.83500549999999
$ws.Range("J11").Value = 200.7586077
$ws.Range("J12").Value = 432.87870320000002
$ws.Range("J13").Value = 227.20696989999999
